$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("D2").Value = "Fail"
$wsTestCases.Range("D3").Value = "Fail"

$wsLogin001 = $wb.Worksheets.Item("Login_001")
$wsLogin001.Range("G8").Value = "Fail"

$wsLogin002 = $wb.Worksheets.Item("Login_002")
$wsLogin002.Range("G8").Value = "Fail"
